# Auto-generated script applying commit 'Add data for 2024-03-23'
# Updates 2024 (column K) violent crime totals, and a handful of
# already-reported 2023 columns (J) / 2018 column (E) that were revised,
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1575
$ws.Range("K3").Value = 1498
$ws.Range("E4").Value = 2027
$ws.Range("J4").Value = 1795
$ws.Range("K4").Value = 326
$ws.Range("K6").Value = 1942
$ws.Range("E7").Value = 26032
$ws.Range("J7").Value = 29264
$ws.Range("K7").Value = 5438

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 100
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 38
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 67
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 39
$ws.Range("K3").Value = 60
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 343
$ws.Range("K11").Value = 112
$ws.Range("K15").Value = 49
$ws.Range("K19").Value = 141
$ws.Range("K20").Value = 125
$ws.Range("K25").Value = 22
$ws.Range("K27").Value = 63
$ws.Range("K29").Value = 248
$ws.Range("K32").Value = 10
$ws.Range("K33").Value = 224
$ws.Range("K36").Value = 60
$ws.Range("K37").Value = 184
$ws.Range("K40").Value = 10
$ws.Range("J42").Value = 1238
$ws.Range("K42").Value = 191
$ws.Range("K52").Value = 146
$ws.Range("K53").Value = 82
$ws.Range("J54").Value = 575
$ws.Range("K54").Value = 92
$ws.Range("E63").Value = 364
$ws.Range("J63").Value = 91
$ws.Range("K63").Value = 20
$ws.Range("K64").Value = 33
$ws.Range("K65").Value = 140
$ws.Range("K66").Value = 22
$ws.Range("K67").Value = 208
$ws.Range("K76").Value = 76
$ws.Range("K82").Value = 8
$ws.Range("K83").Value = 111
$ws.Range("K84").Value = 37
$ws.Range("K85").Value = 280
$ws.Range("K86").Value = 37
$ws.Range("K92").Value = 23
$ws.Range("K95").Value = 92
$ws.Range("E101").Value = 26032
$ws.Range("J101").Value = 29264
$ws.Range("K101").Value = 5438

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J4").Value = 49
$ws.Range("K4").Value = 5
$ws.Range("J7").Value = 575
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 65
$ws.Range("K4").Value = 12
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 45
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 258
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 51
$ws.Range("K6").Value = 83
$ws.Range("J7").Value = 1238
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 91
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 8

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 146
